# Generate Report for Handback
# Refresh the handoff/handback timestamps on the zh-cn and de-de report sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-21 03:45:39"
$wsZhCn.Range("H2").Value = "2016-03-21 03:46:18"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-21 03:45:47"
$wsDeDe.Range("H2").Value = "2016-03-21 03:46:32"
